$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "In Translation"
#    This shared string is referenced by every "Status" cell across all three
#    sheets (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4), so a single
#    find/replace sweep over every sheet updates them all consistently.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# ---------------------------------------------------------------------------
# 2. Narrow the "zh-cn" / "de-de" status columns.
#    Overview: columns E and F (zh-cn / de-de)
#    zh-cn sheet: column C (Status)
#    de-de sheet: column C (Status)
#    The COM ColumnWidth setter only resolves to 1/6-character increments,
#    so 12.5 is the closest input that lands on the narrower target width.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
